$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42: switch to border-bottom style (s=6/7) and add empty A42 cell ---
$ws.Range("A24:E24").Copy()
$ws.Range("A42:E42").PasteSpecial(-4122)

# --- Row 43 (style template row 41) ---
$ws.Range("A41:E41").Copy()
$ws.Range("A43:E43").PasteSpecial(-4122)
$ws.Range("A43").Value = "SCRIPT/D73P31A/us3104.ssb"
$ws.Range("B43").Value = 60
$ws.Range("C43").Value = " The [CS:I]Sky Gift[CR]s I received in\nreturn had all sorts of wonderful items inside."
$ws.Range("D43").Value = " В [CS:I]Небесных Подарках[CR], которых\nдарили мне, было множество чудесных вещей."
$ws.Range("E43").Value = " Â [CS:I]Îåáåòîúö Ðïäàñëàö[CR], ëïóïñúö\näàñéìé íîå, áúìï íîïçåòóâï œôäåòîúö âåþåê."
$ws.Range("A43:E43").RowHeight = 43.2

# --- Row 44 (style template row 41) ---
$ws.Range("A41:E41").Copy()
$ws.Range("A44:E44").PasteSpecial(-4122)
$ws.Range("A44").Value = " SCRIPT/D73P31A/us0203.ssb"
$ws.Range("B44").Value = 63
$ws.Range("C44").Value = " Once I used the items, I made\nit to the summit before I knew it!"
$ws.Range("D44").Value = " Эти предметы очень помогли мне\nдобраться до вершины!"
$ws.Range("E44").Value = " Üóé ðñåäíåóú ïœåîû ðïíïãìé íîå\näïáñàóûòÿ äï âåñšéîú!"
$ws.Range("A44:E44").RowHeight = 57.6

# --- Row 45 (style template row 41) ---
$ws.Range("A41:E41").Copy()
$ws.Range("A45:E45").PasteSpecial(-4122)
$ws.Range("A45").Value = "SCRIPT/D73P31A/us0303.ssb"
$ws.Range("B45").Value = 66
$ws.Range("C45").Value = " [CS:I]Sky Gift[CR]s are such\nwonderful treasures!"
$ws.Range("D45").Value = " [CS:I]Небесные Подарки[CR] это очень\nчудесные сокровища!"
$ws.Range("E45").Value = " [CS:I]Îåáåòîúå Ðïäàñëé[CR] üóï ïœåîû\nœôäåòîúå òïëñïâéþà!"
$ws.Range("A45:E45").RowHeight = 43.2

# --- Row 46 (style template row 24) ---
$ws.Range("A24:E24").Copy()
$ws.Range("A46:E46").PasteSpecial(-4122)
$ws.Range("A46").Value = "SCRIPT/D73P31A/us0403.ssb"
$ws.Range("A46:E46").RowHeight = 43.2

# --- Row 47 (style template row 41) ---
$ws.Range("A41:E41").Copy()
$ws.Range("A47:E47").PasteSpecial(-4122)
$ws.Range("A47").Value = "SCRIPT/P01P04A/us2003.ssb"
$ws.Range("B47").Value = 44
$ws.Range("C47").Value = " I haven't returned in a while…"
$ws.Range("D47").Value = " Давно меня здесь не было..."
$ws.Range("E47").Value = " Äàâîï íåîÿ èäåòû îå áúìï..."
$ws.Range("A47:E47").RowHeight = 43.2

# --- Row 48 (style template row 20) ---
$ws.Range("A20:E20").Copy()
$ws.Range("A48:E48").PasteSpecial(-4122)
$ws.Range("A48").Value = "SCRIPT/P01P04A/us2103.ssb"
$ws.Range("B48").Value = 47
$ws.Range("C48").Value = " I wanted a [CS:I]Sky Gift[CR], so I've\nbeen exploring [CS:P]Sky Peak[CR] for some time."
$ws.Range("D48").Value = " Я хотел найти [CS:I]Небесный Подарок[CR],\nпоэтому я долго исследовал [CS:P]Небесный Пик[CR]."
$ws.Range("E48").Value = " Ÿ öïóåì îàêóé [CS:I]Îåáåòîúê Ðïäàñïë[CR],\nðïüóïíô ÿ äïìãï éòòìåäïâàì [CS:P]Îåáåòîúê Ðéë[CR]."
$ws.Range("A48:E48").RowHeight = 43.2

# --- Row 49 (style template row 41) ---
$ws.Range("A41:E41").Copy()
$ws.Range("A49:E49").PasteSpecial(-4122)
$ws.Range("A49").Value = "SCRIPT/P01P04A/us2203.ssb"
$ws.Range("B49").Value = 18
$ws.Range("C49").Value = " Have you heard of [CS:K]Shaymin[CR]'s\nDelivery Service?"
$ws.Range("D49").Value = " Вы знаете про Службу Доставки\n[CS:K]Шейминов[CR]?"
$ws.Range("E49").Value = " Âú èîàåóå ðñï Òìôçáô Äïòóàâëé\n[CS:K]Šåêíéîïâ[CR]?"
$ws.Range("A49:E49").RowHeight = 43.2

# --- Row 50 (style template row 19) ---
$ws.Range("A19:E19").Copy()
$ws.Range("A50:E50").PasteSpecial(-4122)
$ws.Range("B50").Value = 21
$ws.Range("C50").Value = " It's a service that will deliver\n[CS:I]Sky Gift[CR]s to others."
$ws.Range("D50").Value = " Это служба, которая доставляет\n[CS:I]Небесные Подарки[CR] другим Покемонам."
$ws.Range("E50").Value = " Üóï òìôçáà, ëïóïñàÿ äïòóàâìÿåó\n[CS:I]Îåáåòîúå Ðïäàñëé[CR] äñôãéí Ðïëåíïîàí."
$ws.Range("A50:E50").RowHeight = 31.8

# --- Row 51 (style template row 19) ---
$ws.Range("A19:E19").Copy()
$ws.Range("A51:E51").PasteSpecial(-4122)
$ws.Range("B51").Value = 24
$ws.Range("C51").Value = " The shop is near the [CS:P]Shaymin\nVillage[CR]'s entrance. You should check it\nout sometime."
$ws.Range("D51").Value = " Служба находится у входа в\n[CS:P]Деревню Шейминов[CR]. Вам стоит как-нибудь\nею воспользоваться."
$ws.Range("E51").Value = " Òìôçáà îàöïäéóòÿ ô âöïäà â\n[CS:P]Äåñåâîý Šåêíéîïâ[CR]. Âàí òóïéó ëàë-îéáôäû\nåý âïòðïìûèïâàóûòÿ."
$ws.Range("A51:E51").RowHeight = 42

# --- Restore view state: selection on D50 ---
$ws.Range("D50").Select()

